$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.296.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -5.01%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.564.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -5.02%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'  +0.17%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'288.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.65%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.3755"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.80%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'49.32"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.28%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.3413"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.48%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.68%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07640"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -5.15%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.05%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'21.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.38%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.004"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -4.70%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.925"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.51%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.559.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -5.18%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.00001128"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -6.79%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'89.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -5.65%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06716"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.86%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D21").Value = "'6.229"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -5.92%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'16.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -4.99%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.5270"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -8.35%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'11.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.06%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'22.288.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -5.08%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.395"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.96%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.799"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -7.28%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'20.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -4.29%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'145.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.95%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.972"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.05%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'125.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.85%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.734.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -5.11%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.017"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.75%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'6.187"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -10.02%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.014"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -5.74%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'10.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -9.95%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.08528"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.89%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.02537"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -5.81%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.2315"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -4.32%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'5.495"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -7.20%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.322"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.34%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.06384"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -6.18%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'11.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -9.20%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.6343"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -7.92%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'14.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -9.02%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.10%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.5976"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -6.50%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'3.746"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.49%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -7.05%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.259"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.45%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'124.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.37%  "
$ws.Range("E51").Style = "Normal"
